$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Add the new "keytypes" worksheet at the end of the workbook
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "keytypes"

# ------------------------------------------------------------------
# 2. Update the "molgenis" directory sheet (sheet1) with a new row
#    describing the "keytypes" template
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("molgenis")

$ws1.Range("A10").Value = "keytypes"
$ws1.Range("B10").Value = "template"

$ws1.Range("J10").Value = "http://semanticscience.org/resource/SIO_000762"
$ws1.Hyperlinks.Add($ws1.Range("J10"), "http://semanticscience.org/resource/SIO_000762")
$ws1.Range("J10").Style = $ws1.Range("J9").Style

$ws1.Range("K10").Value = "A database key is an informational entity whose value is constructed from one or more database columns."

# ------------------------------------------------------------------
# 3. Populate the new "keytypes" sheet
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("keytypes")

$ws4.Range("A1").Value = "value"
$ws4.Range("B1").Value = "description"
$ws4.Range("C1").Value = "codesystem"
$ws4.Range("D1").Value = "code"
$ws4.Range("E1").Value = "iri"
$ws4.Range("A1:E1").Style = $ws1.Range("A1").Style

$ws4.Range("A2").Value = "foreign database key"
$ws4.Range("B2").Value = "A foreign database key is a database key that refers to a key in some table."
$ws4.Range("C2").Value = "SIO"
$ws4.Range("D2").Value = 500
$ws4.Range("E2").Value = "http://semanticscience.org/resource/SIO_000764"
$ws4.Hyperlinks.Add($ws4.Range("E2"), "http://semanticscience.org/resource/SIO_000764")

$ws4.Range("A3").Value = "primary database key"
$ws4.Range("B3").Value = "A primary database key is a database key that identifies every row of a table."
$ws4.Range("C3").Value = "SIO"
$ws4.Range("D3").Value = 499
$ws4.Range("E3").Value = "http://semanticscience.org/resource/SIO_000763"
$ws4.Hyperlinks.Add($ws4.Range("E3"), "http://semanticscience.org/resource/SIO_000763")

$ws4.Range("E2").Style = $ws1.Range("J9").Style
$ws4.Range("E3").Style = $ws1.Range("J9").Style

# Put the original active sheet back in focus, matching the source file
$ws1.Select()
